$d = $word.ActiveDocument

# Locate the run containing the sentence that needs to be split/corrected.
# (the trailing "lá" typo -> "là", plus a split into 3 runs per the target
# revision.)
$target = $d.Content
$found = $target.Find.Execute(
    " thời gian lái xe liên tục quá 4h hoặc tổng thời gian làm việc của cùng 1 người lá quá 10h/ ngày.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target sentence to edit."
}

# Re-derive a fresh Range over the same span. Calling InsertXML directly on
# a Range that Find.Execute just repositioned leaves stale state behind, so
# we re-anchor with Start/End first.
$target = $d.Range($target.Start, $target.End)

# Replace the matched range's contents with three runs: the lead-in text,
# the corrected word "là" in its own run, and the trailing text - mirroring
# the exact run split introduced by the edit.
$newXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
  "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
  "<pkg:xmlData>" +
  "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
  "<w:body><w:p>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> thời gian lái xe liên tục quá 4h hoặc tổng thời gian làm việc của cùng 1 người </w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>là</w:t></w:r>" +
  "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> quá 10h/ ngày.</w:t></w:r>" +
  "</w:p></w:body></w:document>" +
  "</pkg:xmlData></pkg:part></pkg:package>"

$target.InsertXML($newXml)

Write-Host "Done."
